$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headlines = 'these analyst revise their forecast on viridian therapeutic following q result kelly clarksons exhusband brandon blackstock dead at after cancer battle arctic pablo coin open final presale stage after raising m and burning over billion token france top court strike down pesticide clause in farming law investor confidence wane a jobys loss mount and commercial timeline remains elusive veteran fund manager turn head with palantir stock price target coinmarketcap listing and early entry bonus analyst say ruvi ai ruvi ha a clear path to gain like early ethereum eth did these analyst revise their forecast on wixcom following q result the year bond auction went poorly yield are popping higher barrons jstar holding announces closing of underwriter overallotment option in connection with initial public offering shiba inu shib holder are jumping to this audited ai token early entry bonus and cmc listing made analyst call it the next millionaire maker eli lilly share plummet toward worst day in yearsheres why it down forbes dogecoin news today dogecoin face bearish indicator a consolidation continues near macys is selling a cotton sheet set for thats cool and comfortable consumer energy the principal subsidiary of cm energy declares quarterly dividend on preferred stock if you bought ozak ai with in july you could be sitting on by the end of august jstar holding announces closing of underwriter overallotment option in connection with initial public offering age verification in gta online how it impact gamers north american energy opportunity corp welcome murray macleod stewart to executive leadership and acquires equity interest in mjolnir resource corp tesla stock get crazy prediction from ceo elon musk draftkings stock surrender q earnings pop whats going on magnite analyst boost their forecast after betterthanexpected q earnings the latest a trump sweeping new tariff take effect u trading partner scramble worldwide the latest a trump sweeping new tariff take effect u trading partner scramble worldwide fastly tap new cfo amid record revenue what u consumer can expect from new tariff on imported good ltl proving to be big growth engine at rxo these analyst slash their forecast on vertex after q result sempra stock rise after second quarter earnings exceed analyst estimate tim scannell former stryker president coo named chair of longeviti neuro solution amid m capital infusion from soleus capital abbvie inc nyseabbv share sold by beacon harbor wealth advisor inc tim scannell former stryker president coo named chair of longeviti neuro solution amid m capital infusion from soleus capital tailoredspace expands with new del mar coworking location state of ohio welcome beta technology aircraft ap technology summarybrief at pm edt booming real estate market in quebec city cma exceptional sale and price in july with legionnaire case growing new york must take a comprehensive approach to address disease two merging boston bank to rebrand under new beacon name crypto news today trump to allow crypto in k plan euphoria raise m paxos settle with ny regulator conocophillips q fy earnings call transcript smb acquiring bb payable fuel priority q topline jump apple win wall street backing a billion u investment eas trump tension boost growth plan dow dip over point u initial claim rise what to know about trump newest and most sweeping tariff dentsply sirona maintains outlook despite sale decline orland library renovation project complete intel stock plummet amid ceo resignation demand trump tariff threat sprout social analyst cut their forecast after q earnings paramount and skydance close their billion merger kicking off reign of new entertainment giant trump is calling for intels ceo to step down here why warren buffett buy stock of company with million howard stern decision to make change in the board of director change in the board of director here what the new tariff could cost each u household here what the new tariff could cost each u household here what the new tariff could cost each u household here what the new tariff could cost each u household here what the new tariff could cost each u household here what the new tariff could cost each u household here what the new tariff could cost each u household here what the new tariff could cost each u household here what the new tariff could cost each u household whoa massive million judgment against assurance iq ftc slam assurance iq for misleading health plan marketing in violation of the ftc act and telemarketing sale rule can opendoor survive the real estate deep freeze these analyst increase their forecast on applovin after strong q earnings crypto group claim responsibility for hurling sex toy at wnba game cces sullivan fresh farmer market feed community support local farm whale activity signal major breakout ahead for ripple xrp and new eth layer token below mgic investment corporation receives moody rating upgrade lyno ai begin first presale round offering aidriven crypto arbitrage access at calcium aluminate cement market size will attain usd million by growing at cagr exclusive report by zion market research global calcium aluminate cement market size share trend analysis report traveler will have a new international flight option outside of big airline calcium aluminate cement market size will attain usd million by growing at cagr exclusive report by zion market research global calcium aluminate cement market size share trend analysis report universal uvv q earnings call transcript datacapable appoints duke energy veteran sandy buzzard to transform utility customer engagement former x ceo linda yaccarino named top exec of weightloss company former x ceo linda yaccarino named top exec of weightloss company datacapable appoints duke energy veteran sandy buzzard to transform utility customer engagement millennials tell u how your finance have changed under the trump administration toyota brace for billion hit from new u tariff amid trade tension developer buy worth of farmland to build luxury home in lehigh valley ramaco resource inc announces closing of upsized million public offering trump push u pickup truck export in trade deal despite tight road abroad realty income o q earnings call transcript new percent additional tariff on import from india ap business summarybrief at pm edt what u consumer can expect from new tariff on imported good what u consumer can expect from new tariff on imported good what u consumer can expect from new tariff on imported good what u consumer can expect from new tariff on imported good united airline resolve tech problem that hit more than flight but residual delay likely check out what whale are doing with axon here why sam altman say openais gpt fall short of agi business insider this is what whale are betting on netflix analyst cheer disney dtc profit park strength but flag capex linear tv risk rxst investor alert robbins geller rudman dowd llp announces that rxsight inc investor with substantial loss have opportunity to lead security class action lawsuit m grant via community west bank fuel visalia affordable housing how do fiduciary advisor protect your financial interest top low price crypto presale project you can still afford cold wallet x entry stuns investor aircraft airframe fuel system market is forecasted to reach u billion in say stratview research gap already emerging in u and russian view of potential trumpputin summit kelly clarksons exhusband brandon blackstock dead at after cancer battle two pennsylvania state trooper have been shot governor say mamdani call alleged cuomotrump conversation disqualifying harvard scientist say research could be set back year after funding freeze south park come for kristi noem and ice american tell u something you actually love about american culture rural voice fourth season to cover immigration health care mental health social isolation and tribal right in rural minnesota top u firm modelled economic scenario for relocating palestinian from gaza outstanding scifi western procedural or action series maybe it time to bring back the genre emmy court vacates fed interchange rule dealing blow to bank stephen miller run the number on why dems are melting down over trump deportation and census change trump seek to change how census collect data and want to exclude immigrant in u illegally letter attack on jewish student must end fbi force out more leader including exdirector who fought trump demand for jan agent name ap news summary at pm edt fbi ousts more leader including exdirector who resisted trump demand for jan agent name letter mac and the tax cut note from a recovering conservative state of ohio welcome beta technology aircraft doe rosalias silence on gaza mean that activism only happens online netanyahu say israel intends to take control of all of gaza trans men wut welp here the dumbest thing youll read maybe ever about masculinity no really a top republican in the georgia governor race is suing his rival over campaign financing here what the new tariff could cost each u household here what the new tariff could cost each u household here what the new tariff could cost each u household here what the new tariff could cost each u household california ag confirms ota may decline to apply conflicting tax regulation in individual appeal a top republican in the georgia governor race is suing his rival over campaign financing a top republican in the georgia governor race is suing his rival over campaign financing what u consumer can expect from new tariff on imported good two die when helicopter strike powerlines and crash into a barge on the mississippi river republican call for investigation into state senator who posted ice sighting rxst investor alert robbins geller rudman dowd llp announces that rxsight inc investor with substantial loss have opportunity to lead security class action lawsuit former superman turn into an ice supervillain will india save russia t armata mbt what the useu trade deal mean for you elon musk plan for a new political party will likely be derailed by a u political system hostile to new voice trump say intel ceo is highly conflicted call for resignation anderson man sentenced for impersonating homeland security officer and possessing incendiary bomb all diva no wow jasmine crockett staffer spill the tea about their horrible noshow bos from hell republican senator claim rfk jr mrna vaccine decision undermines trump agenda doug lamalfa to host town hall in chico on monday republican are out to get janet mill steve collins netanyahu say israel plan to take over all of gaza in bid to destroy hamas this one shocking stat prof that the american dream is dying what they did there wa torture u fisherman disappeared decade ago leaving only his empty boat authority just found his dna in a boot gov maura healey sign law that strengthens transgender reproductive health care protection swiss to seek more talk with u a horror tariff kick in netanyahu say israel plan to take over all of gaza in bid to destroy hamas man who brutally beat elderly prolifers dodge prison time who wear the pant in the family charlie bryant to seek second term on conyers city council ice ha a new courthouse tactic get immigrant case tossed then arrest them outside eleventh circuit expands government authority under the mandatory victim restitution act to collect unsatisfied restitution obligation new leadership appointed in haiti a gang threaten to overthrow government democrat slam texas senator cornyn over alleged fbi role in locating fleeing lawmaker prioritizing protection of civilian in future peace operation applying lesson from year of poc maine is not affordable for young people letter former acting director of the fbi and other official pushed out source say trump want to change the census here how the census work trump want to change the census here how the census work opinion what should voter look for in a new hennepin county attorney trump is honoring purple heart recipient including who sent him medal after attempt on his life trump is honoring purple heart recipient including who sent him medal after attempt on his life bonus editorial carton for aug texas democrat plead for donation to extend their walkout and block trump redistricting plan texas democrat plead for donation to extend their walkout and block trump redistricting plan texas democrat plead for donation to extend their walkout and block trump redistricting plan texas democrat plead for donation to extend their walkout and block trump redistricting plan texas democrat plead for donation to extend their walkout and block trump redistricting plan texas democrat plead for donation to extend their walkout and block trump redistricting plan texas democrat plead for donation to extend their walkout and block trump redistricting plan texas democrat plead for donation to extend their walkout and block trump redistricting plan police seek possible sexual assault victim of unlicensed massage business their deportation proceeding were closed for year trump official are reviving them their deportation proceeding were closed for year trump official are reviving them their deportation proceeding were closed for year trump official are reviving them their deportation proceeding were closed for year trump official are reviving them their deportation proceeding were closed for year trump official are reviving them their deportation proceeding were closed for year trump official are reviving them their deportation proceeding were closed for year trump official are reviving them their deportation proceeding were closed for year trump official are reviving them senior fbi official who resisted trump administration demand reportedly pushed out senior fbi official who resisted trump administration demand reportedly pushed out texas redistricting saga take a shocking turn takeaway from aps report on alaska native response to oil and mining proposal netanyahu say israel plan to take over gaza a official discus widening it offensive trump push for drilling mining sharpens debate for alaska native about land they view a sacred opinion federal official should have to pas fitness test fan point out the irony in former superman actor dean cain joining ice french algerian tie hit new low after macron hardens stance america need big ball not mahmoud khalil vance boelter pleads not guilty in minnesota lawmaker shooting gop strategist lay out how census disproportionately hurt republican netanyahu say israel to control not govern gaza th circuit opts in to the evolving collective action debate russian drone test nato article defense guarantee ahead of friday sanction deadline specialty clothing store facing tariff on good from india concerned about future'

$ws.Range("A33").Value = 45877
$ws.Range("A33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B33").Value = 6389.4501953125
$ws.Range("C33").Value = 6395.16015625
$ws.Range("D33").Value = 6355.22021484375
$ws.Range("E33").Value = 6355.22021484375
$ws.Range("F33").Value = 4769910000
$ws.Range("G33").Value = 0.0077997153489748
$ws.Range("H33").Value = 1
$ws.Range("I33").Value = $headlines
$ws.Range("J33").Value = $headlines
$ws.Range("K33").Value = 0.11
$ws.Range("L33").Value = 0.108
$ws.Range("M33").Value = 0.782
$ws.Range("N33").Value = -0.9875
$ws.Range("O33").Value = 76
$ws.Range("P33").Value = 98
$ws.Range("Q33").Value = 0
$ws.Range("R33").Value = 136
$ws.Range("S33").Value = 76
$ws.Range("T33").Value = 193
$ws.Range("U33").Value = 124
$ws.Range("V33").Value = 49
$ws.Range("W33").Value = 43
$ws.Range("X33").Value = 41

Write-Output "Row 33 written"
